$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: unwrap the error message from brackets and add the Statement ID column
$ws.Range("K2").Value = "The 'Encoded Statement' field does not contain IG Script-encoded content."
$ws.Range("L2").Value = "'1"

# Row 3
$ws.Range("K3").Value = "OK"
$ws.Range("L3").Value = "'2.1"
$ws.Range("AO3").Value = "[OR].Bdir.[2.2]"

# Row 4
$ws.Range("K4").Value = "OK"
$ws.Range("L4").Value = "'2.2"
$ws.Range("AO4").Value = "[OR].Bdir.[2.1]"

# Row 5
$ws.Range("K5").Value = "OK"
$ws.Range("L5").Value = "'3"

# Row 6
$ws.Range("K6").Value = "OK"
$ws.Range("L6").Value = "'4.1"
$ws.Range("AO6").Value = "[OR].Bdir.[4.2]"

# Row 7
$ws.Range("K7").Value = "OK"
$ws.Range("L7").Value = "'4.2"
$ws.Range("AO7").Value = "[OR].Bdir.[4.1]"
